# "Ridimensionato per Lenovo E550"
# Adds a new "dimensioni" worksheet (after the existing sheets) that
# computes screen/window sizing values, and makes it the active sheet.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the last existing sheet so the tab order becomes
# eventi, objMapping, dimensioni.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "dimensioni"

# Labels (entered in this order so shared-string ids line up: home, top,
# bottom, view, backup, piantina)
$ws.Range("A1").Value = "home"
$ws.Range("A2").Value = "top"
$ws.Range("A4").Value = "bottom"
$ws.Range("A3").Value = "view"
$ws.Range("D1").Value = "backup"
$ws.Range("F3").Value = "piantina"

# Numeric data
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 80
$ws.Range("D2").Value = 80

$ws.Range("B3").Value = 6
$ws.Range("C3").Value = 90
$ws.Range("D3").Value = 120
$ws.Range("G3").Value = 9

$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 135
$ws.Range("D4").Value = 135

# Formulas
$ws.Range("C1").FormulaArray = "=SUM(B2:B4*C2:C4)"
$ws.Range("H3").Formula = "=INT(B3*C3/G3)-2"
$ws.Range("I3").Formula = "=H3*G3"

# Selection / active cell on the new sheet
$ws.Range("I7").Select()
